# feat: add 2022-Q4 data
#
# - "总计" gets a new row: the existing 2022-Q1 summary row is turned into
#   the 2022-Q4 summary row (new numbers), and a fresh row is appended below
#   it carrying the original 2022-Q1 numbers.
# - The existing "2022-Q1" sheet (holding per-fund holdings) is duplicated so
#   its data is preserved unchanged in a sheet that stays named "2022-Q1".
#   The original sheet is then overwritten with the new Q4 per-fund holdings
#   and renamed "2022-Q4", ending up positioned before "2022-Q1".

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Item("2022-Q1")

# ---------------------------------------------------------------------
# 1) "总计": shift the old 2022-Q1 row down to row 3, and turn row 2 into
#    the new 2022-Q4 row.
# ---------------------------------------------------------------------
$total.Range("A2").Copy($total.Range("A3"))
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.01

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.02

# ---------------------------------------------------------------------
# 2) Duplicate "2022-Q1" — the copy lands right after the source sheet and
#    keeps the historical per-fund data untouched.
# ---------------------------------------------------------------------
$q1.Copy($null, $q1)
$q1copy = $wb.Worksheets.Item("2022-Q1 (2)")

# Move the original sheet out of the name's way, then rename the copy back
# to "2022-Q1" so the untouched historical data keeps that name.
$q1.Name = "2022-Q4"
$q1copy.Name = "2022-Q1"
$q4 = $q1

# ---------------------------------------------------------------------
# 3) Overwrite the (now named) "2022-Q4" sheet with the new per-fund
#    holdings, re-using the "总计" header's border/bold style for the
#    header row and the A-column marker cell.
# ---------------------------------------------------------------------
$total.Range("B1").Copy($q4.Range("B1:H1"))
$total.Range("A2").Copy($q4.Range("A2"))

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0

# Force plain-text storage (no numeric coercion, no leftover number-format
# style) by writing a text-literal formula, then converting that formula to
# its literal value in place via copy / paste-values.
function Set-TextValue($range, $text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $q4.Range("B2") "004332"
Set-TextValue $q4.Range("C2") "恒生前海沪港深新兴产业精选混合"
Set-TextValue $q4.Range("D2") "0.49"
Set-TextValue $q4.Range("E2") "75.81"
Set-TextValue $q4.Range("F2") "3.20"
Set-TextValue $q4.Range("G2") "0.0157"

$q4.Range("H2").Value = 9
